$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.501.38'
$ws.Range('E2').Value = '  -2.12%  '
$ws.Range('D3').Value = '3.496.27'
$ws.Range('E3').Value = '  -3.30%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '576.99'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -4.64%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '190.19'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -6.73%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.612'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.81%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '3.485.22'
$ws.Range('E8').Value = '  -3.29%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.204'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -5.89%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.618'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.68%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '50.66'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -5.36%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000284'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -6.32%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.10'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -5.56%  '
$ws.Range('D15').Value = '4.058.18'
$ws.Range('E15').Value = '  -3.01%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '638.12'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -6.78%  '
$ws.Range('D17').Value = '69.294.26'
$ws.Range('E17').Value = '  -2.42%  '
$ws.Range('D18').Value = '3.494.54'
$ws.Range('E18').Value = '  -3.32%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.30'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -4.34%  '
$ws.Range('E20').Value = '  -2.22%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '18.33'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.32%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.948'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -4.98%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '18.29'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.40%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.28'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.66%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '98.71'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -9.61%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.27'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -7.81%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.87'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -5.38%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.93'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -6.26%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.31'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -9.49%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '32.36'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -5.83%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.06'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -11.57%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.66'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -8.86%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '587.72'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +12.83%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.53'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -5.65%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.109'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -5.21%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '60.84'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.39%  '
$ws.Range('D37').Value = '3.784.43'
$ws.Range('E37').Value = '  -3.07%  '
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = '0.0₃0788'
$ws.Range('E39').Value = '  -6.71%  '
$ws.Range('B40').Value = 'CoreDAO'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.77'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +33.28%  '
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.83'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -6.21%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.370'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.71%  '
$ws.Range('E44').Value = '  -5.98%  '
$ws.Range('E45').Value = '  -9.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0439'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -5.85%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.32'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.69%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.82'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -8.37%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.135'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.72%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.997'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.14'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.87%  '
